# repull data, push all data, mean calculation
# Update dSF (column F) values for a handful of rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 0
$ws.Range("F8").Value = 5
$ws.Range("F9").Value = 6
$ws.Range("F14").Value = -4
$ws.Range("F17").Value = -3
$ws.Range("F27").Value = -5
